# Trade #8 closed at 2026-02-18 00:09:47 - unknown UNKNOWN +0.000%
#
# Updates the live trading results workbook to reflect that Trade #36
# (row id 36 in the trade logs) on the MarketMaking strategy was closed
# (early exit) instead of remaining OPEN, and rolls that result up into
# the Summary and Strategy Status sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.65   # Current Capital
$summary.Range("B4").Value = 0.75      # Total P&L $
$summary.Range("B5").Value = 0.42      # Total P&L %
$summary.Range("B6").Value = 36        # Total Trades
$summary.Range("B8").Value = 14        # Losing Trades
$summary.Range("B9").Value = 52.78     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 6)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 99.65000000000001   # Capital
$status.Range("D6").Value = 7                   # Trades
$status.Range("E6").Value = -0.16               # P&L $
$status.Range("F6").Value = -0.35               # P&L %
$status.Range("G6").Value = 42.86               # Win Rate %

# ---------------------------------------------------------------------
# All Trades sheet - Trade #36 row (row 37)
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("G37").Value = 0.4             # Exit Price
$allTrades.Range("H37").Value = "CLOSED"        # Status
$allTrades.Range("I37").Value = -29.4815        # P&L %
$allTrades.Range("J37").Value = -0.17           # P&L $
$allTrades.Range("K37").Value = 99.65000000000001  # Capital After
$allTrades.Range("L37").Value = "early_exit"    # Exit Reason
$allTrades.Range("M37").Value = 0.17            # Duration (min)

# ---------------------------------------------------------------------
# MarketMaking sheet - Trade #36 row (row 8)
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
$marketMaking.Range("G8").Value = 0.4           # Exit Price
$marketMaking.Range("H8").Value = "CLOSED"      # Status
$marketMaking.Range("I8").Value = -29.4815      # P&L %
$marketMaking.Range("J8").Value = -0.17         # P&L $
$marketMaking.Range("K8").Value = 99.65000000000001  # Capital After
$marketMaking.Range("P8").Value = "early_exit"  # Exit Reason
$marketMaking.Range("Q8").Value = 0.17          # Duration (min)
